$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1338
$ws.Range("I28").Value = 460.7857
$ws.Range("K28").Value = 460.7857
$ws.Range("M28").Value = 24.21429999999998
$ws.Range("H51").Value = 41675390
$ws.Range("I51").Value = 55564290
$ws.Range("K51").Value = 55564290
$ws.Range("M51").Value = -55563806
$ws.Range("H138").Value = 293761.97
$ws.Range("I138").Value = 3199.7097
$ws.Range("J138").Value = 457533.44
$ws.Range("K138").Value = 9599.1291
$ws.Range("L138").Value = 1372600.32
$ws.Range("M138").Value = -4459.1291
$ws.Range("N138").Value = -1382880.32

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3346.1904
$ws.Range("I61").Value = 2083.6667
$ws.Range("J61").Value = 6502.5
$ws.Range("K61").Value = 2083.6667
$ws.Range("L61").Value = 6502.5
$ws.Range("M61").Value = -1871.6667
$ws.Range("N61").Value = -6926.5
$ws.Range("H74").Value = 133606.88
$ws.Range("I74").Value = 174657.75
$ws.Range("K74").Value = 174657.75
$ws.Range("M74").Value = -173783.75
$ws.Range("H77").Value = 133606.88
$ws.Range("I77").Value = 174657.75
$ws.Range("K77").Value = 873288.75
$ws.Range("M77").Value = -868920.75
$ws.Range("H88").Value = 6507.5713
$ws.Range("I88").Value = 4662.6665
$ws.Range("J88").Value = 7891.25
$ws.Range("K88").Value = 4662.6665
$ws.Range("L88").Value = 7891.25
$ws.Range("M88").Value = -4256.6665
$ws.Range("N88").Value = -8703.25
$ws.Range("H91").Value = 6507.5713
$ws.Range("I91").Value = 4662.6665
$ws.Range("J91").Value = 7891.25
$ws.Range("K91").Value = 4662.6665
$ws.Range("L91").Value = 7891.25
$ws.Range("M91").Value = -3258.6665
$ws.Range("N91").Value = -10699.25
$ws.Range("H122").Value = 5045.05
$ws.Range("I122").Value = 4731.6
$ws.Range("K122").Value = 14194.8
$ws.Range("M122").Value = -11744.8
$ws.Range("H136").Value = 3346.1904
$ws.Range("I136").Value = 2083.6667
$ws.Range("J136").Value = 6502.5
$ws.Range("K136").Value = 6251.000100000001
$ws.Range("L136").Value = 19507.5
$ws.Range("M136").Value = -3701.000100000001
$ws.Range("N136").Value = -24607.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 1658.4546
$ws.Range("I7").Value = 1838.3334
$ws.Range("J7").Value = 849
$ws.Range("K7").Value = 1838.3334
$ws.Range("L7").Value = 849
$ws.Range("M7").Value = -1725.3334
$ws.Range("N7").Value = -1075
$ws.Range("H16").Value = 1496.0817
$ws.Range("I16").Value = 1427.9722
$ws.Range("J16").Value = 1684.6923
$ws.Range("K16").Value = 1427.9722
$ws.Range("L16").Value = 1684.6923
$ws.Range("M16").Value = -1140.9722
$ws.Range("N16").Value = -2258.6923
$ws.Range("H22").Value = 1493.4445
$ws.Range("I22").Value = 1057
$ws.Range("K22").Value = 1057
$ws.Range("M22").Value = -707
$ws.Range("H31").Value = 6368.381
$ws.Range("I31").Value = 5751.909
$ws.Range("J31").Value = 7046.5
$ws.Range("K31").Value = 5751.909
$ws.Range("L31").Value = 7046.5
$ws.Range("M31").Value = -5456.909
$ws.Range("N31").Value = -7636.5
$ws.Range("H34").Value = 6368.381
$ws.Range("I34").Value = 5751.909
$ws.Range("J34").Value = 7046.5
$ws.Range("K34").Value = 5751.909
$ws.Range("L34").Value = 7046.5
$ws.Range("M34").Value = -5549.909
$ws.Range("N34").Value = -7450.5
$ws.Range("H113").Value = 1496.0817
$ws.Range("I113").Value = 1427.9722
$ws.Range("J113").Value = 1684.6923
$ws.Range("K113").Value = 1427.9722
$ws.Range("L113").Value = 1684.6923
$ws.Range("M113").Value = 742.0278000000001
$ws.Range("N113").Value = -6024.6923
$ws.Range("H132").Value = 3039.9285
$ws.Range("I132").Value = 2381.842
$ws.Range("J132").Value = 4429.222
$ws.Range("K132").Value = 7145.526
$ws.Range("L132").Value = 13287.666
$ws.Range("M132").Value = -4615.526
$ws.Range("N132").Value = -18347.666
$ws.Range("H134").Value = 11061.923
$ws.Range("I134").Value = 14575.571
$ws.Range("K134").Value = 43726.713
$ws.Range("M134").Value = -41191.713
$ws.Range("H141").Value = 664962.9
$ws.Range("J141").Value = 664962.9
$ws.Range("L141").Value = 664962.9
$ws.Range("N141").Value = -675322.9

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 791.3333
$ws.Range("I122").Value = 809
$ws.Range("J122").Value = 775.875
$ws.Range("K122").Value = 7281
$ws.Range("L122").Value = 6982.875
$ws.Range("M122").Value = -4831
$ws.Range("N122").Value = -11882.875
$ws.Range("H125").Value = 2999.6667
$ws.Range("I125").Value = 1999.5
$ws.Range("K125").Value = 5998.5
$ws.Range("M125").Value = -1078.5
$ws.Range("H132").Value = 2583.3809
$ws.Range("I132").Value = 1544
$ws.Range("J132").Value = 3528.2727
$ws.Range("K132").Value = 13896
$ws.Range("L132").Value = 31754.4543
$ws.Range("M132").Value = -11366
$ws.Range("N132").Value = -36814.4543
$ws.Range("H133").Value = 6676.6665
$ws.Range("I133").Value = 8515
$ws.Range("J133").Value = 3000
$ws.Range("K133").Value = 25545
$ws.Range("L133").Value = 9000
$ws.Range("M133").Value = -20485
$ws.Range("N133").Value = -19120
$ws.Range("H140").Value = 7211.9
$ws.Range("I140").Value = 6846.5557
$ws.Range("K140").Value = 20539.6671
$ws.Range("M140").Value = -15359.6671

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 2547
$ws.Range("I107").Value = 3866.3333
$ws.Range("K107").Value = 3866.3333
$ws.Range("M107").Value = -1946.3333
$ws.Range("H113").Value = 3179.3333
$ws.Range("I113").Value = 2968.5757
$ws.Range("K113").Value = 2968.5757
$ws.Range("M113").Value = -798.5756999999999
$ws.Range("H122").Value = 7967.968
$ws.Range("I122").Value = 7121.3335
$ws.Range("K122").Value = 21364.0005
$ws.Range("M122").Value = -18914.0005
$ws.Range("H126").Value = 7716.0713
$ws.Range("I126").Value = 6821.3
$ws.Range("K126").Value = 20463.9
$ws.Range("M126").Value = -17993.9
$ws.Range("H132").Value = 1299.5714
$ws.Range("I132").Value = 1299.5714
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 3898.7142
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -1368.7142
$ws.Range("N132").ClearContents()
$ws.Range("H134").Value = 57539.8
$ws.Range("J134").Value = 57539.8
$ws.Range("L134").Value = 172619.4
$ws.Range("N134").Value = -177689.4

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1169.4546
$ws.Range("I16").Value = 1086.4
$ws.Range("J16").Value = 2000
$ws.Range("K16").Value = 1086.4
$ws.Range("L16").Value = 2000
$ws.Range("M16").Value = -916.4000000000001
$ws.Range("N16").Value = -2340
$ws.Range("H40").Value = 83075.92
$ws.Range("I40").Value = 96363.45
$ws.Range("K40").Value = 96363.45
$ws.Range("M40").Value = -96227.45
$ws.Range("H82").Value = 1249.8
$ws.Range("I82").Value = 1249.8
$ws.Range("K82").Value = 1249.8
$ws.Range("M82").Value = -888.8
$ws.Range("H85").Value = 1249.8
$ws.Range("I85").Value = 1249.8
$ws.Range("K85").Value = 1249.8
$ws.Range("M85").Value = -1.799999999999955
$ws.Range("H122").Value = 5183.1665
$ws.Range("I122").Value = 5519.8
$ws.Range("J122").Value = 3500
$ws.Range("K122").Value = 16559.4
$ws.Range("L122").Value = 10500
$ws.Range("M122").Value = -14109.4
$ws.Range("N122").Value = -15400
$ws.Range("H132").Value = 5372.6
$ws.Range("I132").Value = 4106.92
$ws.Range("K132").Value = 12320.76
$ws.Range("M132").Value = -9790.76
$ws.Range("H133").Value = 102314.336
$ws.Range("J133").Value = 102314.336
$ws.Range("L133").Value = 102314.336
$ws.Range("N133").Value = -107374.336
$ws.Range("H136").Value = 3235.9556
$ws.Range("I136").Value = 2933.6765
$ws.Range("K136").Value = 8801.029500000001
$ws.Range("M136").Value = -6251.029500000001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 543.5
$ws.Range("I113").Value = 468.58823
$ws.Range("K113").Value = 1405.76469
$ws.Range("M113").Value = 764.23531
$ws.Range("H122").Value = 9619254
$ws.Range("I122").Value = 4478.8945
$ws.Range("K122").Value = 13436.6835
$ws.Range("M122").Value = -10986.6835
$ws.Range("H132").Value = 3126.5715
$ws.Range("I132").Value = 3121.1155
$ws.Range("K132").Value = 9363.3465
$ws.Range("M132").Value = -6833.3465
